# The deck's slide master currently uses the "Integral" theme (theme1.xml),
# while an unused/auxiliary "Office Theme" definition sits in theme2.xml
# (only wired up to the Notes Master). The authored change swaps the two
# theme color definitions so the deck's actual (slide-facing) theme becomes
# the default Office palette.
#
# Recolor the deck's theme (ThemeColorScheme) to the stock "Office" theme
# color values, in DrawingML clrScheme order:
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# PowerPoint's ColorFormat.RGB uses the VBA RGB() BGR-packed integer, i.e.
# 0x00BBGGRR rather than the usual 0xRRGGBB web-hex ordering.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

$officeThemeColorsBGR = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColorsBGR[$i - 1]
}
